$wb = $excel.ActiveWorkbook

# --- "Create" sheet: update row 2 values ---
$wsCreate = $wb.Worksheets.Item("Create")
$wsCreate.Range("A2").Value = "ActualSMS"
$wsCreate.Range("B2").Value = "Expected"
$wsCreate.Range("D2").Value = "Chat"
$wsCreate.Range("E2").Value = "'404"
$wsCreate.Range("F2").Value = "'404"
$wsCreate.Range("G2").Value = "Actual Test"
$wsCreate.Range("H2").Value = "'404"

# --- "Edit" sheet: update row 2 values ---
$wsEdit = $wb.Worksheets.Item("Edit")
$wsEdit.Range("A2").Value = "ActualSMS"
$wsEdit.Range("B2").Value = "Expected"
$wsEdit.Range("D2").Value = "Chat"
$wsEdit.Range("E2").Value = "'404"
$wsEdit.Range("F2").Value = "'404"
$wsEdit.Range("G2").Value = "Actual Test"
$wsEdit.Range("H2").Value = "'404"
$wsEdit.Range("J2").Value = "Eldorado"

# --- "Delete" sheet: update row 2 values ---
$wsDelete = $wb.Worksheets.Item("Delete")
$wsDelete.Range("A2").Value = "ActualSMS"
$wsDelete.Range("B2").Value = "Expected"
$wsDelete.Range("D2").Value = "Chat"
$wsDelete.Range("E2").Value = "'404"
$wsDelete.Range("F2").Value = "'404"
$wsDelete.Range("G2").Value = "Eldorado"
$wsDelete.Range("H2").Value = "'404"

# --- Update selections on each sheet (no change in active tab for these two) ---
[void]$wsCreate.Range("E2").Select()
[void]$wsEdit.Range("J2").Select()

# --- "Delete" becomes the active sheet/tab, with its own new selection ---
[void]$wsDelete.Activate()
[void]$wsDelete.Range("G2").Select()
